# "formato del excel actualizado"
# Update column F (image/logo file name) on sheet PRODUCTOS: a block of
# rows whose F cell currently holds the placeholder numeric value 0 gets
# replaced with the text "logoproducto.jpeg". Also restore the last
# selected cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRODUCTOS")

$rows = @(2,3,7,8,13,14,17,27,33,35,42,47,48,51,54,55,56,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,77,78,79,80,82,84,86,89,92,93,94,95,96,97,99,100,101,102,105,107,109,110,111,118,120,121,124,125,129,131,132,133,134,135,140,141,143,144,145,149,150,151,152,153,154,155,156,158,160,168,169,171,172,173,174,180,181,183,184,185,186,188,189,191,192,193,195,196,197,198,199,200,203,204,205,206,207,208,209,210,212,214,216,217,220,221,231,232,233,234,263,264,265,266,267,268,269,271,272,273,274,275,276,277,278,279,280,281,282,283,284,289,290,291,292,293,296,297,298,299)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.ClearFormats()
    $cell.Value = "logoproducto.jpeg"
}

$ws.Activate()
$ws.Range("F259").Select()
